# pabi_account_report - xlsx_report_purchase_invoice_plan.xlsx
# #3724 - InvoicePlan Module: add Exchange Rate PO / KV, Acceptance Date,
# PO State, PO Close and InvoicePlan status columns to the report header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Exchange Rate" -> "Exchange Rate PO" (column AK keeps its position/style)
$ws.Range("AK10").Value2 = "Exchange Rate PO"

# A new "Acceptance Date" column is inserted right after "WA Number" (AO),
# pushing every following header one slot to the right.
$ws.Range("AP10").Value2 = "Acceptance Date"
$ws.Range("AQ10").Value2 = "Recieive Quantity"
$ws.Range("AR10").Value2 = "Unit Price"
$ws.Range("AS10").Value2 = "Subtotal"
$ws.Range("AT10").Value2 = "Billing Number"
$ws.Range("AU10").Value2 = "KV Number"
$ws.Range("AV10").Value2 = "Doc Date"
$ws.Range("AW10").Value2 = "Posting Date"
$ws.Range("AX10").Value2 = "Supplier Invoice Number"

# A new "Exchange Rate KV" column is inserted right after "Supplier Invoice
# Number" (now AX), pushing everything after it one more slot to the right.
$ws.Range("AY10").Value2 = "Exchange Rate KV"
$ws.Range("AZ10").Value2 = "Amount "
$ws.Range("BA10").Value2 = "Amount Local"
$ws.Range("BB10").Value2 = "Deposit"
$ws.Range("BC10").Value2 = "Advance/Deposit"
$ws.Range("BD10").Value2 = "FinLease"

# Finally, three brand new trailing columns are appended.
$ws.Range("BE10").Value2 = "PO State "
$ws.Range("BF10").Value2 = "PO Close"
$ws.Range("BG10").Value2 = "InvoicePlan status"

# The author left the whole header row selected when the file was saved.
[void]$ws.Rows(10).Select()
